$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 30, 0.5, 20, 99),
    @(5, 30, 0.5, 20, 99),
    @(6, 30, 0.5, 20, 99),
    @(7, 30, 0.5, 20, 99)
)

$row = 5
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

$ws.Range("B8").Select()
